$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the second "UF" header (column E, destination state) to "UF1"
# to disambiguate it from the first "UF" header in column C (origin state).
$ws.Range("E1").Value = "UF1"

# Column C ("UF" - origin state) had stale/incorrect values in several rows;
# correct them to match column E ("UF1" - destination state) for that row,
# mirroring the values actually present in each row.
$ws.Range("C2").Value = "RJ"
$ws.Range("C3").Value = "MG"
$ws.Range("C4").Value = "MT"
$ws.Range("C5").Value = "PR"
$ws.Range("C6").Value = "PR"
$ws.Range("C7").Value = "PR"
$ws.Range("C8").Value = "PR"
$ws.Range("C12").Value = "SP"
$ws.Range("C13").Value = "SP"
$ws.Range("C14").Value = "BA"
$ws.Range("C15").Value = "ES"
$ws.Range("C16").Value = "PR"
$ws.Range("C17").Value = "PR"
$ws.Range("C18").Value = "SP"
$ws.Range("C20").Value = "MG"
$ws.Range("C21").Value = "MT"
$ws.Range("C22").Value = "MT"
$ws.Range("C23").Value = "PR"
$ws.Range("C24").Value = "PR"
$ws.Range("C25").Value = "PR"
$ws.Range("C26").Value = "PR"
$ws.Range("C27").Value = "PR"
$ws.Range("C28").Value = "PR"
$ws.Range("C29").Value = "PR"
$ws.Range("C33").Value = "SP"

# Underline formatting on two "Destino_Regiões" cells (D6, D9) picked up
# while reviewing/correcting the rows above.
$ws.Range("D6").Font.Underline = $true
$ws.Range("D9").Font.Underline = $true
